$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new blank column before the existing "Late" column (column N),
# pushing Late / Heading / Outstanding one column to the right (O, P, Q).
$ws.Columns("N:N").Insert()

# The newly inserted column keeps a custom (non-autofit) width, matching
# the neighbouring "In Advance" column's width.
$ws.Columns("N").ColumnWidth = 10.33

# Switch the active sheet from "Transactions" to "Repayment Schedule" and
# move the selection there.
$ws.Activate()
$ws.Range("R10").Select()
